$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F7 becomes a (single-space) string instead of the numeric fraction it held before.
# F12 / F17 lose their numeric values entirely (cells become blank).
# A new column K ("jsa") is added holding parsed-JSON / array literal test
# strings for each data row, in the original row order except where noted.
#
# The cells are written in the same order the strings were first introduced
# into the workbook (matches the shared-string table append order).

$ws.Range("F7").Value = ' '
$ws.Range("K1").Value = 'jsa'
$ws.Range("K4").Value = '{"ds":3}'
$ws.Range("K5").Value = '{"ds":4}'
$ws.Range("K6").Value = '[1,3,2]'
$ws.Range("K7").Value = '[1,3,3]'
$ws.Range("K8").Value = '[1,3,4]'
$ws.Range("K9").Value = '["21",3,5]'
$ws.Range("K10").Value = '["21",3,6]'
$ws.Range("K11").Value = '["21",3,"123"]'
$ws.Range("K3").Value = '["21",3,"124"]'
$ws.Range("K13").Value = '["21",3,"125"]'
$ws.Range("K16").Value = '["22","ds","ase"]'
$ws.Range("K17").Value = '{"a":["23","ds","ase"]}'
$ws.Range("K18").Value = '{"a":["24","ds","ase"],"b":1}'
$ws.Range("K2").Value = '{"a":["24","ds","ase"],"b":2}'
$ws.Range("K20").Value = '{"a":["24","ds","ase"],"b":3}'
$ws.Range("K15").Value = '[[12,33,11],"ds","ase"]'
$ws.Range("K14").Value = '["21",{"wq":"re"},"126"]'

# Reuses an already-existing shared string (same text as elsewhere in the sheet).
$ws.Range("K12").Value = '{"d1":4,"dv":{"fd":9}}'

# Remove the now-orphaned numeric values.
$ws.Range("F12").Value = ''
$ws.Range("F17").Value = ''

# Selection moved as part of the edit.
[void]$ws.Range("E8").Select()
